$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 113
$ws.Range("F4").Value = 552
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 1903
$ws.Range("F8").Value = 5002
$ws.Range("F9").Value = 1401
$ws.Range("F11").Value = 2983
$ws.Range("F14").Value = 1241
$ws.Range("F15").Value = 4049
$ws.Range("F16").Value = 952
$ws.Range("F18").Value = 1611
$ws.Range("F19").Value = 68
$ws.Range("F20").Value = 2564
$ws.Range("F21").Value = 17
$ws.Range("F22").Value = 85
$ws.Range("F24").Value = 926
$ws.Range("F25").Value = 273
$ws.Range("F28").Value = 1039
$ws.Range("F29").Value = 309
$ws.Range("F30").Value = 90
$ws.Range("F32").Value = 169
$ws.Range("F33").Value = 1553
$ws.Range("F34").Value = 2091
$ws.Range("F35").Value = 989
$ws.Range("F36").Value = 29
$ws.Range("F37").Value = 231
$ws.Range("F41").Value = 637
$ws.Range("F42").Value = 360
$ws.Range("F43").Value = 251
$ws.Range("F45").Value = 111

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 665

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 665
$ws.Range("F7").Value = 552
$ws.Range("F8").Value = 15
$ws.Range("F9").Value = 1903
$ws.Range("F10").Value = 5002
$ws.Range("F11").Value = 1401
$ws.Range("F14").Value = 2983
$ws.Range("F16").Value = 1241
$ws.Range("F17").Value = 4049
$ws.Range("F18").Value = 952
$ws.Range("F20").Value = 1611
$ws.Range("F22").Value = 68
$ws.Range("F23").Value = 2564
$ws.Range("F24").Value = 27
$ws.Range("F25").Value = 17
$ws.Range("F30").Value = 926
$ws.Range("F31").Value = 273
$ws.Range("F34").Value = 1039
$ws.Range("F35").Value = 309
$ws.Range("F36").Value = 1553
$ws.Range("F37").Value = 2091
$ws.Range("F39").Value = 989
$ws.Range("F40").Value = 29
$ws.Range("F42").Value = 231
$ws.Range("F45").Value = 637
$ws.Range("F46").Value = 360
$ws.Range("F47").Value = 251
$ws.Range("F49").Value = 111
